# Auto-generated Excel COM-interop script applying the Leviathan_Profits value updates
$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$wsALC.Range("H17").Value = 12309.618
$wsALC.Range("J17").Value = 12309.618
$wsALC.Range("L17").Value = 36928.854
$wsALC.Range("N17").Value = -37264.854
$wsALC.Range("H38").Value = 393.33334
$wsALC.Range("I38").Value = 417.72726
$wsALC.Range("J38").Value = 125
$wsALC.Range("K38").Value = 1253.18178
$wsALC.Range("L38").Value = 375
$wsALC.Range("M38").Value = -881.1817799999999
$wsALC.Range("N38").Value = -1119
$wsALC.Range("H39").Value = 1292
$wsALC.Range("I39").Value = 1292
$wsALC.Range("K39").Value = 3876
$wsALC.Range("M39").Value = -3580
$wsALC.Range("H42").Value = 409.42856
$wsALC.Range("I42").Value = 168
$wsALC.Range("K42").Value = 504
$wsALC.Range("M42").Value = -274
$wsALC.Range("H55").Value = 530.8333
$wsALC.Range("J55").Value = 620
$wsALC.Range("L55").Value = 620
$wsALC.Range("N55").Value = -1048
$wsALC.Range("H95").Value = 27961.75
$wsALC.Range("J95").Value = 27961.75
$wsALC.Range("L95").Value = 27961.75
$wsALC.Range("N95").Value = -33453.75
$wsALC.Range("H103").Value = 17857590
$wsALC.Range("I103").Value = 362.5
$wsALC.Range("J103").Value = 20833794
$wsALC.Range("K103").Value = 1087.5
$wsALC.Range("L103").Value = 62501382
$wsALC.Range("M103").Value = -501.5
$wsALC.Range("N103").Value = -62502554
$wsALC.Range("H112").Value = 2068.7827
$wsALC.Range("J112").Value = 2232.4736
$wsALC.Range("L112").Value = 6697.4208
$wsALC.Range("N112").Value = -8913.4208
$wsALC.Range("H123").Value = 68634.75
$wsALC.Range("J123").Value = 68634.75
$wsALC.Range("L123").Value = 68634.75
$wsALC.Range("N123").Value = -78434.75
$wsALC.Range("H132").Value = 2337.0535
$wsALC.Range("I132").Value = 765.55817
$wsALC.Range("J132").Value = 7535.077
$wsALC.Range("K132").Value = 2296.67451
$wsALC.Range("L132").Value = 22605.231
$wsALC.Range("M132").Value = 233.3254900000002
$wsALC.Range("N132").Value = -27665.231

# --- ARM ---
$wsARM.Range("H19").Value = 0
$wsARM.Range("J19").Value = 0
$wsARM.Range("L19").Value = 0
$wsARM.Range("N19").ClearContents()
$wsARM.Range("H61").Value = 893.6
$wsARM.Range("I61").Value = 783.2727
$wsARM.Range("J61").Value = 1197
$wsARM.Range("K61").Value = 783.2727
$wsARM.Range("L61").Value = 1197
$wsARM.Range("M61").Value = -571.2727
$wsARM.Range("N61").Value = -1621
$wsARM.Range("H74").Value = 2083.7273
$wsARM.Range("I74").Value = 1614.7778
$wsARM.Range("J74").Value = 2408.3845
$wsARM.Range("K74").Value = 1614.7778
$wsARM.Range("L74").Value = 2408.3845
$wsARM.Range("M74").Value = -740.7778000000001
$wsARM.Range("N74").Value = -4156.3845
$wsARM.Range("H77").Value = 2083.7273
$wsARM.Range("I77").Value = 1614.7778
$wsARM.Range("J77").Value = 2408.3845
$wsARM.Range("K77").Value = 8073.889
$wsARM.Range("L77").Value = 12041.9225
$wsARM.Range("M77").Value = -3705.889
$wsARM.Range("N77").Value = -20777.9225
$wsARM.Range("H95").Value = 30208
$wsARM.Range("I95").Value = 0
$wsARM.Range("J95").Value = 30208
$wsARM.Range("K95").Value = 0
$wsARM.Range("L95").Value = 30208
$wsARM.Range("M95").ClearContents()
$wsARM.Range("N95").Value = -35700
$wsARM.Range("H136").Value = 893.6
$wsARM.Range("I136").Value = 783.2727
$wsARM.Range("J136").Value = 1197
$wsARM.Range("K136").Value = 2349.8181
$wsARM.Range("L136").Value = 3591
$wsARM.Range("M136").Value = 200.1819
$wsARM.Range("N136").Value = -8691

# --- BSM ---
$wsBSM.Range("H134").Value = 1375.9259
$wsBSM.Range("J134").Value = 5000
$wsBSM.Range("L134").Value = 15000
$wsBSM.Range("N134").Value = -20070

# --- CRP ---
$wsCRP.Range("H31").Value = 1351.1482
$wsCRP.Range("I31").Value = 1233.1364
$wsCRP.Range("K31").Value = 1233.1364
$wsCRP.Range("M31").Value = -938.1364000000001
$wsCRP.Range("H32").Value = 5851
$wsCRP.Range("J32").Value = 4500
$wsCRP.Range("L32").Value = 4500
$wsCRP.Range("N32").Value = -5132
$wsCRP.Range("H34").Value = 1351.1482
$wsCRP.Range("I34").Value = 1233.1364
$wsCRP.Range("K34").Value = 1233.1364
$wsCRP.Range("M34").Value = -1031.1364
$wsCRP.Range("H43").Value = 16552.666
$wsCRP.Range("J43").Value = 16552.666
$wsCRP.Range("L43").Value = 16552.666
$wsCRP.Range("N43").Value = -16920.666
$wsCRP.Range("H101").Value = 16552.666
$wsCRP.Range("J101").Value = 16552.666
$wsCRP.Range("L101").Value = 16552.666
$wsCRP.Range("N101").Value = -23042.666
$wsCRP.Range("H122").Value = 28402.764
$wsCRP.Range("I122").Value = 39123.703
$wsCRP.Range("K122").Value = 117371.109
$wsCRP.Range("M122").Value = -114921.109

# --- CUL ---
$wsCUL.Range("H12").Value = 700.0909
$wsCUL.Range("I12").Value = 800
$wsCUL.Range("J12").Value = 690.1
$wsCUL.Range("K12").Value = 2400
$wsCUL.Range("L12").Value = 2070.3
$wsCUL.Range("M12").Value = -2227
$wsCUL.Range("N12").Value = -2416.3
$wsCUL.Range("H68").Value = 2425.4443
$wsCUL.Range("J68").Value = 2478.625
$wsCUL.Range("L68").Value = 7435.875
$wsCUL.Range("N68").Value = -9057.875
$wsCUL.Range("H71").Value = 2425.4443
$wsCUL.Range("J71").Value = 2478.625
$wsCUL.Range("L71").Value = 22307.625
$wsCUL.Range("N71").Value = -30419.625
$wsCUL.Range("H120").Value = 20000
$wsCUL.Range("J120").Value = 20000
$wsCUL.Range("L120").Value = 60000
$wsCUL.Range("N120").Value = -69676
$wsCUL.Range("H126").Value = 5751.5
$wsCUL.Range("I126").Value = 1500
$wsCUL.Range("J126").Value = 10003
$wsCUL.Range("K126").Value = 4500
$wsCUL.Range("L126").Value = 30009
$wsCUL.Range("M126").Value = 440
$wsCUL.Range("N126").Value = -39889

# --- GSM ---
$wsGSM.Range("H70").Value = 5845.5293
$wsGSM.Range("I70").Value = 5789.6665
$wsGSM.Range("J70").Value = 5979.6
$wsGSM.Range("K70").Value = 5789.6665
$wsGSM.Range("L70").Value = 5979.6
$wsGSM.Range("M70").Value = -5519.6665
$wsGSM.Range("N70").Value = -6519.6
$wsGSM.Range("H73").Value = 5845.5293
$wsGSM.Range("I73").Value = 5789.6665
$wsGSM.Range("J73").Value = 5979.6
$wsGSM.Range("K73").Value = 5789.6665
$wsGSM.Range("L73").Value = 5979.6
$wsGSM.Range("M73").Value = -4853.6665
$wsGSM.Range("N73").Value = -7851.6
$wsGSM.Range("H97").Value = 835.85
$wsGSM.Range("I97").Value = 672.0714
$wsGSM.Range("J97").Value = 1218
$wsGSM.Range("K97").Value = 672.0714
$wsGSM.Range("L97").Value = 1218
$wsGSM.Range("M97").Value = -176.0714
$wsGSM.Range("N97").Value = -2210
$wsGSM.Range("H101").Value = 26885.334
$wsGSM.Range("J101").Value = 26885.334
$wsGSM.Range("L101").Value = 26885.334
$wsGSM.Range("N101").Value = -33375.334
$wsGSM.Range("H102").Value = 3483.5833
$wsGSM.Range("I102").Value = 3435.3635
$wsGSM.Range("K102").Value = 3435.3635
$wsGSM.Range("M102").Value = -1813.3635

# --- LTW ---
$wsLTW.Range("H32").Value = 1448.1111
$wsLTW.Range("I32").Value = 790.5
$wsLTW.Range("J32").Value = 3749.75
$wsLTW.Range("K32").Value = 790.5
$wsLTW.Range("L32").Value = 3749.75
$wsLTW.Range("M32").Value = -473.5
$wsLTW.Range("N32").Value = -4383.75
$wsLTW.Range("H55").Value = 450.05
$wsLTW.Range("I55").Value = 454.64285
$wsLTW.Range("K55").Value = 454.64285
$wsLTW.Range("M55").Value = -281.64285
$wsLTW.Range("H61").Value = 500910
$wsLTW.Range("I61").Value = 500910
$wsLTW.Range("J61").Value = 0
$wsLTW.Range("K61").Value = 500910
$wsLTW.Range("L61").Value = 0
$wsLTW.Range("M61").Value = -500708
$wsLTW.Range("N61").ClearContents()
$wsLTW.Range("H113").Value = 500910
$wsLTW.Range("I113").Value = 500910
$wsLTW.Range("J113").Value = 0
$wsLTW.Range("K113").Value = 500910
$wsLTW.Range("L113").Value = 0
$wsLTW.Range("M113").Value = -498740
$wsLTW.Range("N113").ClearContents()

# --- WVR ---
$wsWVR.Range("H101").Value = 17750.5
$wsWVR.Range("J101").Value = 17750.5
$wsWVR.Range("L101").Value = 17750.5
$wsWVR.Range("N101").Value = -24240.5
